# Add a set of "common Excel error" formula examples to column L of the
# "Final" sheet, one per error type (#DIV/0!, #NAME?, #N/A, #NULL!, #NUM!,
# #REF!, #VALUE!), and leave the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final")

# L3: #DIV/0!  - dividing by zero
$ws.Range("L3").Formula = "=H2/0"

# L4: #NAME?   - reference to an undefined name, entered as an (legacy)
#               array formula, i.e. Ctrl+Shift+Enter
$ws.Range("L4").FormulaArray = "=SUM(A2:A3) + UNKNOWN"

# L5: #N/A     - VLOOKUP can't find the lookup value
$ws.Range("L5").Formula = "=VLOOKUP(""NonexistentValue"", A2:A21, 2, FALSE)"

# L6: #NULL!   - intersection operator (space) between non-intersecting ranges
$ws.Range("L6").Formula = "=SUM(A2:A3 B4:B5)"

# L7: #NUM!    - invalid numeric argument (square root of a negative number)
$ws.Range("L7").Formula = "=SQRT(-1)"

# L8: #REF!    - invalid cell reference
$ws.Range("L8").Formula = "=#REF!+A22"

# L9: #VALUE!  - incompatible operand types
$ws.Range("L9").Formula = "=A2 + ""text"""

# Match the author's final cursor position in the saved workbook.
$ws.Range("E9").Select() | Out-Null
